# Fixed n on leaf size graphs
# Rename the B/C header cells from the old "LeafSize.Script.*" labels to the
# corrected "MaxLeafSize.Script.*" labels, widen columns B:C to fit the new
# (longer) header text, and leave the active selection on H11 (matching the
# state the workbook was saved in after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "MaxLeafSize.Script.LeafPosition"
$ws.Range("C1").Value = "MaxLeafSize.Script.MaxLeafSize"

$ws.Columns("B:C").ColumnWidth = 29

$ws.Range("H11").Select()
